# Update "想去人数" (F) and "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets, reflecting refreshed scrape counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1466
    $ws.Range("F4").Value = 1755
    $ws.Range("F7").Value = 654
    $ws.Range("F12").Value = 79
    $ws.Range("F13").Value = 147
    $ws.Range("F18").Value = 4700
    $ws.Range("G18").Value = 19.9

    if ($sheetName -eq "展览") {
        $ws.Range("F19").Value = 44
        $ws.Range("F20").Value = 821
        $ws.Range("F21").Value = 104
        $ws.Range("F22").Value = 2199
        $ws.Range("F25").Value = 2064
    } else {
        $ws.Range("F20").Value = 44
        $ws.Range("F22").Value = 821
        $ws.Range("F23").Value = 104
        $ws.Range("F24").Value = 2199
        $ws.Range("F27").Value = 2064
    }
}
